# Refresh market-data-driven profit columns (H:N) on the Leve profit
# sheets, mirroring the scheduled runner's nightly price pull.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 252.2963
$ws.Range("I28").Value = 220.26086
$ws.Range("J28").Value = 436.5
$ws.Range("K28").Value = 220.26086
$ws.Range("L28").Value = 436.5
$ws.Range("M28").Value = 264.73914
$ws.Range("N28").Value = -1406.5
$ws.Range("H38").Value = 793.35297
$ws.Range("I38").Value = 73.916664
$ws.Range("J38").Value = 2520
$ws.Range("K38").Value = 221.749992
$ws.Range("L38").Value = 7560
$ws.Range("M38").Value = 150.250008
$ws.Range("N38").Value = -8304
$ws.Range("H99").Value = 651.7778
$ws.Range("I99").Value = 633.25
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 1899.75
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -401.75
$ws.Range("N99").Value = -5396
$ws.Range("H103").Value = 545.2222
$ws.Range("I103").Value = 525.875
$ws.Range("K103").Value = 1577.625
$ws.Range("M103").Value = -991.625
$ws.Range("H111").Value = 3286.25
$ws.Range("I111").Value = 2950
$ws.Range("J111").Value = 3622.5
$ws.Range("K111").Value = 8850
$ws.Range("L111").Value = 10867.5
$ws.Range("M111").Value = -5783
$ws.Range("N111").Value = -17001.5
$ws.Range("H129").Value = 949.5574
$ws.Range("I129").Value = 592.7778
$ws.Range("J129").Value = 1011.3077
$ws.Range("K129").Value = 1778.3334
$ws.Range("L129").Value = 3033.9231
$ws.Range("M129").Value = 3221.6666
$ws.Range("N129").Value = -13033.9231
$ws.Range("H132").Value = 2153.2654
$ws.Range("I132").Value = 1550.9678
$ws.Range("K132").Value = 4652.903399999999
$ws.Range("M132").Value = -2122.903399999999
$ws.Range("H137").Value = 1024.5358
$ws.Range("I137").Value = 768.3333
$ws.Range("J137").Value = 1612.2941
$ws.Range("K137").Value = 2304.9999
$ws.Range("L137").Value = 4836.8823
$ws.Range("M137").Value = 245.0001000000002
$ws.Range("N137").Value = -9936.882300000001
$ws.Range("H138").Value = 2816.5881
$ws.Range("I138").Value = 1900.8695
$ws.Range("J138").Value = 4731.273
$ws.Range("K138").Value = 5702.6085
$ws.Range("L138").Value = 14193.819
$ws.Range("M138").Value = -562.6085000000003
$ws.Range("N138").Value = -24473.819

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1504.1666
$ws.Range("I45").Value = 1137
$ws.Range("J45").Value = 1871.3334
$ws.Range("K45").Value = 1137
$ws.Range("L45").Value = 1871.3334
$ws.Range("M45").Value = -760
$ws.Range("N45").Value = -2625.3334
$ws.Range("H74").Value = 1105.826
$ws.Range("I74").Value = 1061.7
$ws.Range("J74").Value = 1400
$ws.Range("K74").Value = 1061.7
$ws.Range("L74").Value = 1400
$ws.Range("M74").Value = -187.7
$ws.Range("N74").Value = -3148
$ws.Range("H77").Value = 1105.826
$ws.Range("I77").Value = 1061.7
$ws.Range("J77").Value = 1400
$ws.Range("K77").Value = 5308.5
$ws.Range("L77").Value = 7000
$ws.Range("M77").Value = -940.5
$ws.Range("N77").Value = -15736
$ws.Range("H97").Value = 616.7308
$ws.Range("I97").Value = 574
$ws.Range("K97").Value = 574
$ws.Range("M97").Value = -78
$ws.Range("H114").Value = 31750
$ws.Range("J114").Value = 31750
$ws.Range("L114").Value = 31750
$ws.Range("N114").Value = -40428
$ws.Range("H122").Value = 1596.7778
$ws.Range("I122").Value = 1336.5714
$ws.Range("K122").Value = 4009.7142
$ws.Range("M122").Value = -1559.7142

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2437.5527
$ws.Range("I105").Value = 2150.577
$ws.Range("J105").Value = 3059.3333
$ws.Range("K105").Value = 2150.577
$ws.Range("L105").Value = 3059.3333
$ws.Range("M105").Value = -403.5770000000002
$ws.Range("N105").Value = -6553.3333
$ws.Range("H107").Value = 24723.545
$ws.Range("I107").Value = 29583.055
$ws.Range("J107").Value = 2855.75
$ws.Range("K107").Value = 29583.055
$ws.Range("L107").Value = 2855.75
$ws.Range("M107").Value = -27663.055
$ws.Range("N107").Value = -6695.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1556.75
$ws.Range("I16").Value = 1520.25
$ws.Range("J16").Value = 1575
$ws.Range("K16").Value = 1520.25
$ws.Range("L16").Value = 1575
$ws.Range("M16").Value = -1233.25
$ws.Range("N16").Value = -2149
$ws.Range("H92").Value = 48888.777
$ws.Range("J92").Value = 48888.777
$ws.Range("L92").Value = 48888.777
$ws.Range("N92").Value = -53880.777
$ws.Range("H113").Value = 1556.75
$ws.Range("I113").Value = 1520.25
$ws.Range("J113").Value = 1575
$ws.Range("K113").Value = 1520.25
$ws.Range("L113").Value = 1575
$ws.Range("M113").Value = 649.75
$ws.Range("N113").Value = -5915
$ws.Range("H132").Value = 230318.08
$ws.Range("I132").Value = 301185.94
$ws.Range("J132").Value = 2528.5715
$ws.Range("K132").Value = 903557.8200000001
$ws.Range("L132").Value = 7585.7145
$ws.Range("M132").Value = -901027.8200000001
$ws.Range("N132").Value = -12645.7145
$ws.Range("H134").Value = 1618.6666
$ws.Range("I134").Value = 1205.579
$ws.Range("J134").Value = 3188.4
$ws.Range("K134").Value = 3616.737
$ws.Range("L134").Value = 9565.200000000001
$ws.Range("M134").Value = -1081.737
$ws.Range("N134").Value = -14635.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 50003370
$ws.Range("J22").Value = 4490
$ws.Range("L22").Value = 13470
$ws.Range("N22").Value = -13808
$ws.Range("H27").Value = 50003370
$ws.Range("J27").Value = 4490
$ws.Range("L27").Value = 13470
$ws.Range("N27").Value = -13674
$ws.Range("H33").Value = 889.0625
$ws.Range("I33").Value = 350.81818
$ws.Range("K33").Value = 2104.90908
$ws.Range("M33").Value = -1821.90908
$ws.Range("H44").Value = 578
$ws.Range("I44").Value = 470
$ws.Range("J44").Value = 740
$ws.Range("K44").Value = 1410
$ws.Range("L44").Value = 2220
$ws.Range("M44").Value = -1012
$ws.Range("N44").Value = -3016
$ws.Range("H68").Value = 1005.4545
$ws.Range("I68").Value = 868.5714
$ws.Range("J68").Value = 1069.3334
$ws.Range("K68").Value = 2605.7142
$ws.Range("L68").Value = 3208.0002
$ws.Range("M68").Value = -1794.7142
$ws.Range("N68").Value = -4830.0002
$ws.Range("H71").Value = 1005.4545
$ws.Range("I71").Value = 868.5714
$ws.Range("J71").Value = 1069.3334
$ws.Range("K71").Value = 7817.1426
$ws.Range("L71").Value = 9624.000599999999
$ws.Range("M71").Value = -3761.1426
$ws.Range("N71").Value = -17736.0006
$ws.Range("H112").Value = 4747.273
$ws.Range("I112").Value = 3000
$ws.Range("J112").Value = 4801.875
$ws.Range("K112").Value = 9000
$ws.Range("L112").Value = 14405.625
$ws.Range("M112").Value = -7892
$ws.Range("N112").Value = -16621.625
$ws.Range("H122").Value = 862.0952
$ws.Range("J122").Value = 1084.9231
$ws.Range("L122").Value = 9764.3079
$ws.Range("N122").Value = -14664.3079
$ws.Range("H123").Value = 800
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H131").Value = 13528670
$ws.Range("I131").Value = 14000
$ws.Range("J131").Value = 15166811
$ws.Range("K131").Value = 42000
$ws.Range("L131").Value = 45500433
$ws.Range("M131").Value = -36960
$ws.Range("N131").Value = -45510513
$ws.Range("H132").Value = 2087.2942
$ws.Range("I132").Value = 679.8333
$ws.Range("J132").Value = 2855
$ws.Range("K132").Value = 6118.4997
$ws.Range("L132").Value = 25695
$ws.Range("M132").Value = -3588.4997
$ws.Range("N132").Value = -30755

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 101554.664
$ws.Range("I97").Value = 40101.54
$ws.Range("J97").Value = 501000
$ws.Range("K97").Value = 40101.54
$ws.Range("L97").Value = 501000
$ws.Range("M97").Value = -39605.54
$ws.Range("N97").Value = -501992
$ws.Range("H113").Value = 1780.6666
$ws.Range("J113").Value = 2277.4546
$ws.Range("L113").Value = 2277.4546
$ws.Range("N113").Value = -6617.4546
$ws.Range("H122").Value = 2439.4614
$ws.Range("I122").Value = 1873.5
$ws.Range("K122").Value = 5620.5
$ws.Range("M122").Value = -3170.5
$ws.Range("H132").Value = 1474.0385
$ws.Range("I132").Value = 842.6842
$ws.Range("J132").Value = 3187.7144
$ws.Range("K132").Value = 2528.0526
$ws.Range("L132").Value = 9563.143199999999
$ws.Range("M132").Value = 1.947400000000016
$ws.Range("N132").Value = -14623.1432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 27136.625
$ws.Range("I61").Value = 35248.832
$ws.Range("K61").Value = 35248.832
$ws.Range("M61").Value = -35046.832
$ws.Range("H93").Value = 859.82355
$ws.Range("I93").Value = 565.7273
$ws.Range("J93").Value = 1399
$ws.Range("K93").Value = 565.7273
$ws.Range("L93").Value = 1399
$ws.Range("M93").Value = 682.2727
$ws.Range("N93").Value = -3895
$ws.Range("H113").Value = 27136.625
$ws.Range("I113").Value = 35248.832
$ws.Range("K113").Value = 35248.832
$ws.Range("M113").Value = -33078.832
$ws.Range("H136").Value = 2617.9753
$ws.Range("I136").Value = 2654.3691
$ws.Range("J136").Value = 2470.125
$ws.Range("K136").Value = 7963.1073
$ws.Range("L136").Value = 7410.375
$ws.Range("M136").Value = -5413.1073
$ws.Range("N136").Value = -12510.375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 736.75
$ws.Range("I107").Value = 689
$ws.Range("J107").Value = 880
$ws.Range("K107").Value = 2067
$ws.Range("L107").Value = 2640
$ws.Range("M107").Value = -147
$ws.Range("N107").Value = -6480
$ws.Range("H122").Value = 13022429
$ws.Range("I122").Value = 12501330
$ws.Range("J122").Value = 15627929
$ws.Range("K122").Value = 37503990
$ws.Range("L122").Value = 46883787
$ws.Range("M122").Value = -37501540
$ws.Range("N122").Value = -46888687
$ws.Range("H132").Value = 1253.9487
$ws.Range("I132").Value = 924.6786
$ws.Range("J132").Value = 2092.0908
$ws.Range("K132").Value = 2774.0358
$ws.Range("L132").Value = 6276.2724
$ws.Range("M132").Value = -244.0357999999997
$ws.Range("N132").Value = -11336.2724

Write-Host "Updated 262 cells across 8 sheets."
